$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.928.87"
$ws.Range("E2").Value = "  -0.94%  "

$ws.Range("D3").Value = "'1.894.54"
$ws.Range("E3").Value = "  -1.00%  "

$ws.Range("E4").Value = "  -0.44%  "

$ws.Range("D5").Value = "'0.7512"
$ws.Range("E5").Value = "  +2.05%  "

$ws.Range("D6").Value = "'239.72"
$ws.Range("E6").Value = "  -1.73%  "

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("D8").Value = "'0.3045"
$ws.Range("E8").Value = "  -2.34%  "

$ws.Range("D9").Value = "'24.91"
$ws.Range("E9").Value = "  -7.29%  "

$ws.Range("D10").Value = "'0.06822"
$ws.Range("E10").Value = "  -1.63%  "

$ws.Range("D11").Value = "'0.07967"
$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").Value = "'0.7463"
$ws.Range("E12").Value = "  -2.91%  "

$ws.Range("D13").Value = "'1.900.28"
$ws.Range("E13").Value = "  -2.44%  "

$ws.Range("D14").Value = "'5.187"
$ws.Range("E14").Value = "  -1.13%  "

$ws.Range("D15").Value = "'91.06"
$ws.Range("E15").Value = "  +0.13%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "'29.940.82"
$ws.Range("E16").Value = "  -1.04%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "'6.085"
$ws.Range("E17").Value = "  +5.18%  "

$ws.Range("D18").Value = "'13.85"
$ws.Range("E18").Value = "  -1.82%  "

$ws.Range("D19").Value = "'0.000007655"
$ws.Range("E19").Value = "  -1.90%  "

$ws.Range("D20").Value = "'234.41"
$ws.Range("E20").Value = "  -4.07%  "

$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("D22").Value = "'2.155.31"
$ws.Range("E22").Value = "  -1.00%  "

$ws.Range("E23").Value = "  -0.47%  "

$ws.Range("D24").Value = "'6.915"
$ws.Range("E24").Value = "  +4.37%  "

$ws.Range("D25").Value = "'9.209"
$ws.Range("E25").Value = "  -1.46%  "

$ws.Range("D26").Value = "'164.93"
$ws.Range("E26").Value = "  -0.26%  "

$ws.Range("D27").Value = "'18.63"

$ws.Range("D28").Value = "'0.1276"
$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("D29").Value = "'2.037"
$ws.Range("E29").Value = "  -3.81%  "

$ws.Range("D30").Value = "'1.335"
$ws.Range("E30").Value = "  -1.41%  "

$ws.Range("D31").Value = "'1.512"
$ws.Range("E31").Value = "  -2.29%  "

$ws.Range("D32").Value = "'4.264"
$ws.Range("E32").Value = "  -1.50%  "

$ws.Range("D33").Value = "'3.985"
$ws.Range("E33").Value = "  -1.54%  "

$ws.Range("D34").Value = "'0.05352"
$ws.Range("E34").Value = "  +3.98%  "

$ws.Range("D35").Value = "'1.238"
$ws.Range("E35").Value = "  -3.57%  "

$ws.Range("D36").Value = "'0.7273"
$ws.Range("E36").Value = "  -2.29%  "

$ws.Range("D37").Value = "'2.717"
$ws.Range("E37").Value = "  -2.29%  "

$ws.Range("D38").Value = "'0.01922"
$ws.Range("E38").Value = "  -0.32%  "

$ws.Range("D39").Value = "'2.764"
$ws.Range("E39").Value = "  -1.10%  "

$ws.Range("D40").Value = "'6.192"
$ws.Range("E40").Value = "  -2.15%  "

$ws.Range("D41").Value = "'0.4404"
$ws.Range("E41").Value = "  -0.84%  "

$ws.Range("D42").Value = "'72.32"
$ws.Range("E42").Value = "  -4.08%  "

$ws.Range("D43").Value = "'1.913"
$ws.Range("E43").Value = "  -0.91%  "

$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  -0.28%  "

$ws.Range("D45").Value = "'0.8234"
$ws.Range("E45").Value = "  -1.30%  "

$ws.Range("D46").Value = "'101.09"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("D47").Value = "'7.548"
$ws.Range("E47").Value = "  -0.67%  "

$ws.Range("D48").Value = "'9.784"
$ws.Range("E48").Value = "  -0.23%  "

$ws.Range("D49").Value = "'2.057.63"
$ws.Range("E49").Value = "  -0.67%  "

$ws.Range("D50").Value = "'36.06"
$ws.Range("E50").Value = "  -2.45%  "

$ws.Range("D51").Value = "'0.05951"
$ws.Range("E51").Value = "  -0.55%  "
